$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = <new price text>; E = <new volume text> }
# Only columns present for a given row were changed in the source update.
$updates = @{
    2 = @{ D='51.951.28'; E='  -0.47%  ' }
    3 = @{ D='2.816.80'; E='  +0.75%  ' }
    4 = @{ E='  -0.05%  ' }
    5 = @{ D='354.85'; E='  +2.59%  ' }
    6 = @{ D='111.60'; E='  -4.23%  ' }
    7 = @{ D='0.569'; E='  +2.84%  ' }
    8 = @{ D='0.999'; E='  -0.01%  ' }
    9 = @{ D='0.597'; E='  +2.34%  ' }
    10 = @{ D='40.54'; E='  -6.01%  ' }
    11 = @{ D='0.0858'; E='  +0.23%  ' }
    13 = @{ E='  -1.32%  ' }
    14 = @{ E='  -0.50%  ' }
    15 = @{ D='3.259.11'; E='  +0.73%  ' }
    16 = @{ D='2.813.53'; E='  +0.70%  ' }
    17 = @{ D='0.930'; E='  +4.25%  ' }
    18 = @{ D='51.726.34'; E='  -0.81%  ' }
    19 = @{ D='7.51'; E='  +5.78%  ' }
    20 = @{ D='3.16'; E='  -1.65%  ' }
    21 = @{ D='13.41'; E='  +0.12%  ' }
    22 = @{ D='0.0₃0993'; E='  +0.97%  ' }
    23 = @{ D='70.68'; E='  +0.61%  ' }
    24 = @{ D='268.73'; E='  -0.50%  ' }
    25 = @{ E='  +1.57%  ' }
    26 = @{ D='26.93'; E='  +1.01%  ' }
    27 = @{ E='  +0.14%  ' }
    28 = @{ E='  +0.47%  ' }
    29 = @{ E='  +0.94%  ' }
    30 = @{ D='0.0493'; E='  +19.59%  ' }
    31 = @{ E='  +2.60%  ' }
    32 = @{ D='52.51'; E='  +4.49%  ' }
    33 = @{ D='34.56'; E='  -1.24%  ' }
    34 = @{ D='5.94'; E='  +3.64%  ' }
    35 = @{ D='5.55'; E='  +11.43%  ' }
    36 = @{ D='0.0851'; E='  +3.38%  ' }
    37 = @{ E='  -0.06%  ' }
    38 = @{ D='3.27'; E='  +1.16%  ' }
    39 = @{ D='2.04'; E='  -4.13%  ' }
    40 = @{ D='18.33'; E='  -3.46%  ' }
    41 = @{ E='  +1.22%  ' }
    42 = @{ D='126.92' }
    43 = @{ E='  -2.82%  ' }
    44 = @{ D='2.48'; E='  -8.26%  ' }
    45 = @{ D='2.28'; E='  -1.13%  ' }
    46 = @{ D='2.074.31'; E='  -0.14%  ' }
    47 = @{ E='  -0.14%  ' }
    48 = @{ E='  -4.62%  ' }
    49 = @{ D='5.83'; E='  +5.27%  ' }
    50 = @{ E='  -0.85%  ' }
    51 = @{ D='9.13'; E='  +1.77%  ' }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey('D')) {
        # Price column is stored as text in the sheet (e.g. "51.951.28", "111.60")
        # Force text format first so Excel does not reinterpret it as a number
        # and strip significant trailing zeros / re-split on the dot separators.
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals.D
        $cell.ClearFormats()
    }
    if ($vals.ContainsKey('E')) {
        $ws.Range("E$row").Value = $vals.E
    }
}
